# Update ABSM1_RN, M1_RN, CM2_RN, CMN3_RN, CMN4_RN (columns D-H) for several
# rows: recalibrated extrapolation values after removing the < $5 price
# point that was polluting the calibration with noise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @(110409.3908774742, -0.006190087236016062, 0.1670491715926389, -1.375772603870159, 17.01743023179016)
    3  = @(111185.158075578,  -0.01116453934005763,  0.2025668836000296, -1.687299940982553, 19.12339510114475)
    4  = @(111847.7858363563, -0.01607793054692739,  0.2237846441944246, -1.659983810261343, 18.26025916565363)
    5  = @(112537.8142592119, -0.02366714428898253,  0.2389743831241476, -1.427074832248852, 12.97648910173117)
    6  = @(113311.1599209056, -0.01966298012433302,  0.2094600744501642, -0.3772497099170516, 5.410188786629753)
    7  = @(113657.5315285968, -0.03366348536694454,  0.263501596236572,  -1.78366541778223,   15.10594639157191)
    8  = @(115113.3739530644, -0.04850368502288669,  0.2220867862683029, -1.069415029573028,  7.350973012894964)
    9  = @(116843.4880410809, -0.09638973781563547,  0.3759359773370344, -1.491749604097103,  9.312567148337838)
    10 = @(118169.8471346056, -0.1363673245867893,   0.4764656808826944, -1.632750298015933,  7.803622810733291)
    17 = @(110439.2903022914, 0.01386554118168657,   0.1406463871966,    -0.5656738459692482, 9.750533922987451)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
